$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove trailing/leading spaces from word-pair cells (rows 33-41, columns B and C)
$ws.Range("B33").Value = "museum"
$ws.Range("B34").Value = "trip"
$ws.Range("C34").Value = "tur"
$ws.Range("B38").Value = "number"
$ws.Range("B39").Value = "project"
$ws.Range("C40").Value = "succes"
$ws.Range("B41").Value = "mood"

# Update the view/selection to match the saved workbook state
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I31").Select()
